# Burn down chart update: fill in "Actual tasks rem." values for Sprint 2
# rows 10-19 (D10:D19), matching the continuing burn-down trend that was
# already present for rows 2-9. The chart series on the Sprint 2 sheet
# references 'Sprint 2'!$D$2:$D$21, so updating these cells also extends
# the chart's cached data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint 2")

$ws.Range("D10").Value = 20
$ws.Range("D11").Value = 19
$ws.Range("D12").Value = 17
$ws.Range("D13").Value = 16
$ws.Range("D14").Value = 15
$ws.Range("D15").Value = 12
$ws.Range("D16").Value = 10
$ws.Range("D17").Value = 7
$ws.Range("D18").Value = 5
$ws.Range("D19").Value = 3
